# Update the alternative text (alt text / "descr") of the comic-strip
# image in the document. The diff shows the wp:docPr and pic:cNvPr
# "descr" attributes changing from the auto-generated Office description
# ("Diagram / Description automatically generated") to a manually
# authored description of the comic. Word's object model exposes both
# of those XML attributes through the single InlineShape.AlternativeText
# property, so setting it once updates both places in the OOXML.

$d = $word.ActiveDocument

$newAltText = "Comic strip. Cat looks at himself in the mirror and says &quot;It's important to reserve some time for self-reflection. Wouldn't you agree, stranger?&quot;"

if ($d.InlineShapes.Count -ge 1) {
    $shape = $d.InlineShapes.Item(1)
    $shape.AlternativeText = $newAltText
}
